$d = $word.ActiveDocument
$d.Content.Find.Execute("GSoC", $false, $false, $false, $false, $false, $true, 1, $false, "Google Summer of Code", 2) | Out-Null
